$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "X" in the sophistication column (Basic=C, Medium=D, Intensive=E)
# for each category/visualization-method row.
$ws.Range("D8").Value = "X"   # Hierarchical -> Medium
$ws.Range("C9").Value = "X"   # K-means -> Basic
$ws.Range("D10").Value = "X"  # Gaussian mixture model -> Medium

$ws.Range("D12").Value = "X"  # Network Visualization -> Medium

$ws.Range("D14").Value = "X"  # Linear regression -> Medium
$ws.Range("C15").Value = "X"  # Pearson correlation -> Basic
$ws.Range("C16").Value = "X"  # Kernel Density Estimate -> Basic

$ws.Range("D18").Value = "X"  # t-test -> Medium
$ws.Range("D19").Value = "X"  # p-value -> Medium

$ws.Range("E21").Value = "X"  # Boxplots -> Intensive
$ws.Range("E22").Value = "X"  # Violin plots -> Intensive
$ws.Range("E23").Value = "X"  # Histogram -> Intensive

$ws.Range("E25").Value = "X"  # Cartogram map -> Intensive

$ws.Range("C27").Value = "X"  # WordCloud -> Basic
$ws.Range("C28").Value = "X"  # Barplot -> Basic

$ws.Range("D30").Value = "X"  # Line graphs -> Medium
$ws.Range("D31").Value = "X"  # scatter plots -> Medium

# Move the active selection to E31
$ws.Range("E31").Select() | Out-Null
